# Coffee menu update:
#  - Hot Coffee!F2  -> add "Best selling" tag (Espresso)
#  - Hot Coffee!F6  -> add "30% off" tag (Kaapicino)
#  - Hot Coffee row 11 (Flat white):
#       B11 description shortened (drop the quoted tagline)
#       C11 price re-entered as text "230" (was numeric 230)
#       D11 gets a new text price "300"
#       F11 gets a new "Cooling Special" tag

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hot Coffee")

# New "Best selling" tag on the Espresso row
$ws.Cells.Item(2, 6).Value = "Best selling"

# New "30% off" tag on the Kaapicino row
$ws.Cells.Item(6, 6).Value = "30% off"

# Flat white row updates
$ws.Cells.Item(11, 2).Value = "Just The Way Our Australian Friends Like It "

# Re-enter the 250ml price as text (leading apostrophe forces text, like typing it in Excel)
$ws.Cells.Item(11, 3).Value = "'230"

# New 350ml price, also stored as text
$ws.Cells.Item(11, 4).Value = "'300"

# New tag for the Flat white
$ws.Cells.Item(11, 6).Value = "Cooling Special"

# Move the active selection (cosmetic, matches the saved view state)
$ws.Range("J7").Select()
